$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (single decimal point) must be
# forced to Text format first, otherwise Excel auto-converts them to
# numbers and silently drops significant trailing/representation detail
# (e.g. "0.140" -> 0.14, "9.00" -> 9). The source workbook stores every
# cell in this table as an inline string, so we must preserve the exact
# literal text.
$numericLookingCells = @(
    "D5", "D6", "D8", "D9", "D10", "D13", "D14", "D15",
    "D16", "D19", "D20", "D22", "D25", "D26", "D27", "D28",
    "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36",
    "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D46",
    "D47", "D48", "D49", "D50"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '37.867.52'
$ws.Range("E2").Value = '  +1.65%  '
$ws.Range("D3").Value = '2.104.91'
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '234.14'
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '57.85'
$ws.Range("E8").Value = '  +1.08%  '
$ws.Range("D9").Value = '0.391'
$ws.Range("E9").Value = '  +2.06%  '
$ws.Range("D10").Value = '0.0782'
$ws.Range("E10").Value = '  +3.14%  '
$ws.Range("E11").Value = '  +3.63%  '
$ws.Range("D12").Value = '2.401.51'
$ws.Range("E12").Value = '  +1.81%  '
$ws.Range("D13").Value = '14.63'
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = '21.38'
$ws.Range("E14").Value = '  +2.22%  '
$ws.Range("D15").Value = '0.780'
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '5.27'
$ws.Range("E16").Value = '  +2.39%  '
$ws.Range("D17").Value = '2.099.84'
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("D18").Value = '37.814.11'
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("D19").Value = '6.21'
$ws.Range("E19").Value = '  -2.25%  '
$ws.Range("D20").Value = '70.88'
$ws.Range("E20").Value = '  +2.32%  '
$ws.Range("D21").Value = '0.0₃0825'
$ws.Range("E21").Value = '  +1.98%  '
$ws.Range("D22").Value = '228.03'
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("D25").Value = '2.40'
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("D26").Value = '168.11'
$ws.Range("E26").Value = '  +1.30%  '
$ws.Range("D27").Value = '0.140'
$ws.Range("E27").Value = '  +10.35%  '
$ws.Range("D28").Value = '9.00'
$ws.Range("E28").Value = '  +2.31%  '
$ws.Range("D29").Value = '1.43'
$ws.Range("E29").Value = '  -2.30%  '
$ws.Range("D30").Value = '19.53'
$ws.Range("E30").Value = '  +2.72%  '
$ws.Range("D31").Value = '0.119'
$ws.Range("E31").Value = '  +1.34%  '
$ws.Range("D32").Value = '4.69'
$ws.Range("E32").Value = '  +5.22%  '
$ws.Range("D33").Value = '0.0631'
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.66'
$ws.Range("E34").Value = '  +0.87%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '2.58'
$ws.Range("E35").Value = '  +3.65%  '
$ws.Range("D36").Value = '3.46'
$ws.Range("E36").Value = '  +5.60%  '
$ws.Range("D37").Value = '1.82'
$ws.Range("E37").Value = '  +4.52%  '
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").Value = '5.43'
$ws.Range("E39").Value = '  -5.00%  '
$ws.Range("D40").Value = '0.0988'
$ws.Range("E40").Value = '  +6.12%  '
$ws.Range("D41").Value = '2.96'
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Value = '97.84'
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("D43").Value = '0.0216'
$ws.Range("E43").Value = '  +2.32%  '
$ws.Range("D44").Value = '1.458.34'
$ws.Range("E44").Value = '  -1.16%  '
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '1.07'
$ws.Range("E46").Value = '  +4.53%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '4.10'
$ws.Range("E47").Value = '  -6.43%  '
$ws.Range("D48").Value = '15.76'
$ws.Range("E48").Value = '  +4.36%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '7.33'
$ws.Range("E49").Value = '  +2.41%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '3.03'
$ws.Range("E50").Value = '  +2.79%  '
$ws.Range("D51").Value = '2.299.11'
$ws.Range("E51").Value = '  +2.34%  '

# Restore default styling on the cells we temporarily forced to Text
# format, so the resulting cell style matches the rest of the sheet.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
